$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 14.05541930061673
$ws.Cells.Item(2, 3).Value = 10.47391468169997
$ws.Cells.Item(2, 5).Value = 15.651076959552
$ws.Cells.Item(2, 6).Value = 41.10382332000375
$ws.Cells.Item(2, 7).Value = 3.682484561007291
$ws.Cells.Item(2, 10).Value = 8.766960101684132
$ws.Cells.Item(2, 11).Value = 9.770976534452185
$ws.Cells.Item(2, 12).Value = 12.01197037873515
$ws.Cells.Item(2, 13).Value = 16.19327610793412
$ws.Cells.Item(2, 14).Value = 21.37722336015483
$ws.Cells.Item(2, 15).Value = 27.77270199753175
$ws.Cells.Item(3, 2).Value = 13.88278255972271
$ws.Cells.Item(3, 3).Value = 10.48265985621742
$ws.Cells.Item(3, 5).Value = 15.66617811906081
$ws.Cells.Item(3, 6).Value = 41.14498196276925
$ws.Cells.Item(3, 7).Value = 3.684183835158776
$ws.Cells.Item(3, 10).Value = 8.75985112007584
$ws.Cells.Item(3, 11).Value = 9.646502559874977
$ws.Cells.Item(3, 12).Value = 12.01073920887318
$ws.Cells.Item(3, 13).Value = 16.16492527225335
$ws.Cells.Item(3, 14).Value = 21.43783914519221
$ws.Cells.Item(3, 15).Value = 27.85347616530236
$ws.Cells.Item(4, 2).Value = 13.77833394807025
$ws.Cells.Item(4, 3).Value = 10.48841822446777
$ws.Cells.Item(4, 5).Value = 15.67740707764562
$ws.Cells.Item(4, 6).Value = 41.1783620181906
$ws.Cells.Item(4, 7).Value = 3.685283120711587
$ws.Cells.Item(4, 10).Value = 8.75544391008472
$ws.Cells.Item(4, 11).Value = 9.570881667991456
$ws.Cells.Item(4, 12).Value = 12.01144468199336
$ws.Cells.Item(4, 13).Value = 16.14971356544147
$ws.Cells.Item(4, 14).Value = 21.4768063450194
$ws.Cells.Item(4, 15).Value = 27.90770131790789
$ws.Cells.Item(5, 2).Value = 13.73620871804913
$ws.Cells.Item(5, 3).Value = 10.49086284734826
$ws.Cells.Item(5, 5).Value = 15.68247580181723
$ws.Cells.Item(5, 6).Value = 41.19400341533866
$ws.Cells.Item(5, 7).Value = 3.685745194099769
$ws.Cells.Item(5, 10).Value = 8.753637449972508
$ws.Cells.Item(5, 11).Value = 9.54030290958832
$ws.Cells.Item(5, 12).Value = 12.0121009708495
$ws.Cells.Item(5, 13).Value = 16.14407159014682
$ws.Cells.Item(5, 14).Value = 21.49312683861315
$ws.Cells.Item(5, 15).Value = 27.93096192076417
$ws.Cells.Item(6, 2).Value = 13.72924173634151
$ws.Cells.Item(6, 3).Value = 10.4912747046659
$ws.Cells.Item(6, 5).Value = 15.68334724693898
$ws.Cells.Item(6, 6).Value = 41.19672377473415
$ws.Cells.Item(6, 7).Value = 3.685822774241945
$ws.Cells.Item(6, 10).Value = 8.753336859410323
$ws.Cells.Item(6, 11).Value = 9.535240669617071
$ws.Cells.Item(6, 12).Value = 12.0122322583158
$ws.Cells.Item(6, 13).Value = 16.14316851336994
$ws.Cells.Item(6, 14).Value = 21.49586351966198
$ws.Cells.Item(6, 15).Value = 27.93489457204969
$ws.Cells.Item(7, 2).Value = 13.77776399388073
$ws.Cells.Item(7, 3).Value = 10.48845079617064
$ws.Cells.Item(7, 5).Value = 15.67747343978488
$ws.Cells.Item(7, 6).Value = 41.1785647094493
$ws.Cells.Item(7, 7).Value = 3.68528929522366
$ws.Cells.Item(7, 10).Value = 8.755419589806721
$ws.Cells.Item(7, 11).Value = 9.570468265671078
$ws.Cells.Item(7, 12).Value = 12.01145203800071
$ws.Cells.Item(7, 13).Value = 16.1496352148115
$ws.Cells.Item(7, 14).Value = 21.4770246614442
$ws.Cells.Item(7, 15).Value = 27.90801030889037
$ws.Cells.Item(8, 2).Value = 13.99560146252717
$ws.Cells.Item(8, 3).Value = 10.47684948837095
$ws.Cells.Item(8, 5).Value = 15.65587819750469
$ws.Cells.Item(8, 6).Value = 41.11633197498496
$ws.Cells.Item(8, 7).Value = 3.68305888954079
$ws.Cells.Item(8, 10).Value = 8.76451778155068
$ws.Cells.Item(8, 11).Value = 9.727911288002675
$ws.Cells.Item(8, 12).Value = 12.01124337397761
$ws.Cells.Item(8, 13).Value = 16.18304795165703
$ws.Cells.Item(8, 14).Value = 21.39776147900949
$ws.Cells.Item(8, 15).Value = 27.79959166584656
$ws.Cells.Item(9, 2).Value = 14.43293614478399
$ws.Cells.Item(9, 3).Value = 10.45717207239147
$ws.Cells.Item(9, 5).Value = 15.62902010442381
$ws.Cells.Item(9, 6).Value = 41.05861661592858
$ws.Cells.Item(9, 7).Value = 3.679126844282809
$ws.Cells.Item(9, 10).Value = 8.782018258075745
$ws.Cells.Item(9, 11).Value = 10.04153168226222
$ws.Cells.Item(9, 12).Value = 12.02236885153038
$ws.Cells.Item(9, 13).Value = 16.26577099305333
$ws.Cells.Item(9, 14).Value = 21.25614373244112
$ws.Cells.Item(9, 15).Value = 27.623746790924
$ws.Cells.Item(10, 2).Value = 14.75754725892218
$ws.Cells.Item(10, 3).Value = 10.44457169862088
$ws.Cells.Item(10, 5).Value = 15.61867953344081
$ws.Cells.Item(10, 6).Value = 41.05537254866283
$ws.Cells.Item(10, 7).Value = 3.67650454183714
$ws.Cells.Item(10, 10).Value = 8.794661130255827
$ws.Cells.Item(10, 11).Value = 10.27288356915588
$ws.Cells.Item(10, 12).Value = 12.03748806864015
$ws.Cells.Item(10, 13).Value = 16.3367174216822
$ws.Cells.Item(10, 14).Value = 21.16043594238356
$ws.Cells.Item(10, 15).Value = 27.51700655319309
$ws.Cells.Item(11, 2).Value = 14.90533042799121
$ws.Cells.Item(11, 3).Value = 10.4392391220681
$ws.Cells.Item(11, 5).Value = 15.61600231582698
$ws.Cells.Item(11, 6).Value = 41.06237395251232
$ws.Cells.Item(11, 7).Value = 3.675368888363971
$ws.Cells.Item(11, 10).Value = 8.800364617511343
$ws.Cells.Item(11, 11).Value = 10.37790956072235
$ws.Cells.Item(11, 12).Value = 12.0458536467451
$ws.Cells.Item(11, 13).Value = 16.3711269230446
$ws.Cells.Item(11, 14).Value = 21.11868863414902
$ws.Cells.Item(11, 15).Value = 27.47333116217211
$ws.Cells.Item(12, 2).Value = 14.96125893868467
$ws.Cells.Item(12, 3).Value = 10.43727697385146
$ws.Cells.Item(12, 5).Value = 15.61527875733419
$ws.Cells.Item(12, 6).Value = 41.06624064707317
$ws.Cells.Item(12, 7).Value = 3.674947033600104
$ws.Cells.Item(12, 10).Value = 8.802517351384017
$ws.Cells.Item(12, 11).Value = 10.41761431294285
$ws.Cells.Item(12, 12).Value = 12.04923344327505
$ws.Cells.Item(12, 13).Value = 16.38445717218798
$ws.Cells.Item(12, 14).Value = 21.10313617758171
$ws.Cells.Item(12, 15).Value = 27.45749498035494
$ws.Cells.Item(13, 2).Value = 14.94921602579321
$ws.Cells.Item(13, 3).Value = 10.43769701831476
$ws.Cells.Item(13, 5).Value = 15.61542169989597
$ws.Cells.Item(13, 6).Value = 41.06535389470222
$ws.Cells.Item(13, 7).Value = 3.675037523825996
$ws.Cells.Item(13, 10).Value = 8.802054038155386
$ws.Cells.Item(13, 11).Value = 10.4090666820053
$ws.Cells.Item(13, 12).Value = 12.04849615099649
$ws.Cells.Item(13, 13).Value = 16.38157303113438
$ws.Cells.Item(13, 14).Value = 21.10647429463835
$ws.Cells.Item(13, 15).Value = 27.46087431840742
$ws.Cells.Item(14, 2).Value = 14.90993266178852
$ws.Cells.Item(14, 3).Value = 10.43907655034898
$ws.Cells.Item(14, 5).Value = 15.61593697846532
$ws.Cells.Item(14, 6).Value = 41.06266773238635
$ws.Cells.Item(14, 7).Value = 3.675334018150704
$ws.Cells.Item(14, 10).Value = 8.800541866597644
$ws.Cells.Item(14, 11).Value = 10.38117762120403
$ws.Cells.Item(14, 12).Value = 12.04612746907814
$ws.Cells.Item(14, 13).Value = 16.37221763693773
$ws.Cells.Item(14, 14).Value = 21.11740399323787
$ws.Cells.Item(14, 15).Value = 27.47201422149278
$ws.Cells.Item(15, 2).Value = 14.88586459203826
$ws.Cells.Item(15, 3).Value = 10.43992899306188
$ws.Cells.Item(15, 5).Value = 15.6162903629068
$ws.Cells.Item(15, 6).Value = 41.0611805435727
$ws.Cells.Item(15, 7).Value = 3.675516695241439
$ws.Cells.Item(15, 10).Value = 8.799614694053973
$ws.Cells.Item(15, 11).Value = 10.36408510503722
$ws.Cells.Item(15, 12).Value = 12.04470412150122
$ws.Cells.Item(15, 13).Value = 16.36652605960061
$ws.Cells.Item(15, 14).Value = 21.12413209871186
$ws.Cells.Item(15, 15).Value = 27.47892927477136
$ws.Cells.Item(16, 2).Value = 14.7478872771102
$ws.Cells.Item(16, 3).Value = 10.44492821050172
$ws.Cells.Item(16, 5).Value = 15.61889517894799
$ws.Cells.Item(16, 6).Value = 41.05508528324673
$ws.Cells.Item(16, 7).Value = 3.676579908019802
$ws.Cells.Item(16, 10).Value = 8.794287416681314
$ws.Cells.Item(16, 11).Value = 10.26601248410977
$ws.Cells.Item(16, 12).Value = 12.0369711273831
$ws.Cells.Item(16, 13).Value = 16.33451106421322
$ws.Cells.Item(16, 14).Value = 21.16320016100264
$ws.Cells.Item(16, 15).Value = 27.51995913775745
$ws.Cells.Item(17, 2).Value = 14.66323612653867
$ws.Cells.Item(17, 3).Value = 10.44809718495183
$ws.Cells.Item(17, 5).Value = 15.62101142535434
$ws.Cells.Item(17, 6).Value = 41.05351493151409
$ws.Cells.Item(17, 7).Value = 3.677246788449445
$ws.Cells.Item(17, 10).Value = 8.791006932306974
$ws.Cells.Item(17, 11).Value = 10.20576724585263
$ws.Cells.Item(17, 12).Value = 12.03260684964676
$ws.Cells.Item(17, 13).Value = 16.31541301595301
$ws.Cells.Item(17, 14).Value = 21.18762495767739
$ws.Cells.Item(17, 15).Value = 27.54638044205155
$ws.Cells.Item(18, 2).Value = 14.6145608925844
$ws.Cells.Item(18, 3).Value = 10.44995750315599
$ws.Cells.Item(18, 5).Value = 15.62241947992296
$ws.Cells.Item(18, 6).Value = 41.05340987950331
$ws.Cells.Item(18, 7).Value = 3.677635750933334
$ws.Cells.Item(18, 10).Value = 8.789115594293927
$ws.Cells.Item(18, 11).Value = 10.1710973867239
$ws.Cells.Item(18, 12).Value = 12.03023679597976
$ws.Cells.Item(18, 13).Value = 16.30462977579056
$ws.Cells.Item(18, 14).Value = 21.20184206675936
$ws.Cells.Item(18, 15).Value = 27.56203668080914
$ws.Cells.Item(19, 2).Value = 14.59808423166329
$ws.Cells.Item(19, 3).Value = 10.45059384178233
$ws.Cells.Item(19, 5).Value = 15.62292903603437
$ws.Cells.Item(19, 6).Value = 41.05351149431445
$ws.Cells.Item(19, 7).Value = 3.677768373970542
$ws.Cells.Item(19, 10).Value = 8.788474453885227
$ws.Cells.Item(19, 11).Value = 10.15935669318525
$ws.Cells.Item(19, 12).Value = 12.02945846893956
$ws.Cells.Item(19, 13).Value = 16.30101356314762
$ws.Cells.Item(19, 14).Value = 21.20668473010917
$ws.Cells.Item(19, 15).Value = 27.56741649139133
$ws.Cells.Item(20, 2).Value = 14.6722463195937
$ws.Cells.Item(20, 3).Value = 10.44775595175551
$ws.Cells.Item(20, 5).Value = 15.62076640392424
$ws.Cells.Item(20, 6).Value = 41.05359951378006
$ws.Cells.Item(20, 7).Value = 3.677175240278186
$ws.Cells.Item(20, 10).Value = 8.791356611452784
$ws.Cells.Item(20, 11).Value = 10.21218261482198
$ws.Cells.Item(20, 12).Value = 12.03305694223244
$ws.Cells.Item(20, 13).Value = 16.31742523584463
$ws.Cells.Item(20, 14).Value = 21.18500745401465
$ws.Cells.Item(20, 15).Value = 27.54352029646608
$ws.Cells.Item(21, 2).Value = 14.92147245073683
$ws.Cells.Item(21, 3).Value = 10.43866979840615
$ws.Cells.Item(21, 5).Value = 15.61577776147529
$ws.Cells.Item(21, 6).Value = 41.06342376980272
$ws.Cells.Item(21, 7).Value = 3.675246708547654
$ws.Cells.Item(21, 10).Value = 8.800986220878752
$ws.Cells.Item(21, 11).Value = 10.38937138905091
$ws.Cells.Item(21, 12).Value = 12.0468174731217
$ws.Cells.Item(21, 13).Value = 16.37495745504254
$ws.Cells.Item(21, 14).Value = 21.11418672785151
$ws.Cells.Item(21, 15).Value = 27.46872308347982
$ws.Cells.Item(22, 2).Value = 15.08413677750205
$ws.Cells.Item(22, 3).Value = 10.43306466264545
$ws.Cells.Item(22, 5).Value = 15.61420855012945
$ws.Cells.Item(22, 6).Value = 41.07692664427137
$ws.Cells.Item(22, 7).Value = 3.674034036127388
$ws.Cells.Item(22, 10).Value = 8.807238533060996
$ws.Cells.Item(22, 11).Value = 10.50477239934965
$ws.Cells.Item(22, 12).Value = 12.05704507441099
$ws.Cells.Item(22, 13).Value = 16.41430394980584
$ws.Cells.Item(22, 14).Value = 21.06939492680934
$ws.Cells.Item(22, 15).Value = 27.4239352332914
$ws.Cells.Item(23, 2).Value = 14.99735618987088
$ws.Cells.Item(23, 3).Value = 10.43602582410344
$ws.Cells.Item(23, 5).Value = 15.61489174627123
$ws.Cells.Item(23, 6).Value = 41.06907327437535
$ws.Cells.Item(23, 7).Value = 3.67467690737466
$ws.Cells.Item(23, 10).Value = 8.803905384830529
$ws.Cells.Item(23, 11).Value = 10.44322884261699
$ws.Cells.Item(23, 12).Value = 12.05147416547529
$ws.Cells.Item(23, 13).Value = 16.39314662348405
$ws.Cells.Item(23, 14).Value = 21.09316486312519
$ws.Cells.Item(23, 15).Value = 27.44746426989298
$ws.Cells.Item(24, 2).Value = 14.66817283286476
$ws.Cells.Item(24, 3).Value = 10.4479101035575
$ws.Cells.Item(24, 5).Value = 15.62087658185176
$ws.Cells.Item(24, 6).Value = 41.05355878890548
$ws.Cells.Item(24, 7).Value = 3.677207569871586
$ws.Cells.Item(24, 10).Value = 8.791198538140405
$ws.Cells.Item(24, 11).Value = 10.2092823299386
$ws.Cells.Item(24, 12).Value = 12.03285302217156
$ws.Cells.Item(24, 13).Value = 16.31651489868411
$ws.Cells.Item(24, 14).Value = 21.18619028230697
$ws.Cells.Item(24, 15).Value = 27.54481191590701
$ws.Cells.Item(25, 2).Value = 14.31383732070055
$ws.Cells.Item(25, 3).Value = 10.462168102726
$ws.Cells.Item(25, 5).Value = 15.63463255972512
$ws.Cells.Item(25, 6).Value = 41.06734660431608
$ws.Cells.Item(25, 7).Value = 3.68014355688809
$ws.Cells.Item(25, 10).Value = 8.777321118742009
$ws.Cells.Item(25, 11).Value = 9.956378575657904
$ws.Cells.Item(25, 12).Value = 12.01813272319781
$ws.Cells.Item(25, 13).Value = 16.24158204909106
$ws.Cells.Item(25, 14).Value = 21.29298464402561
$ws.Cells.Item(25, 15).Value = 27.66737729327744
